# forests-scraped.xlsx update - 2025-10-16 12:18
# Move the 5 rows currently on "New" into "Previously added" (appended at the
# bottom), then replace "New" with 4 freshly scraped rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

$xlPasteFormats = -4122
$xlShiftUp      = -4162

# ---------------------------------------------------------------------
# 1. Capture the data currently sitting on "New" (rows 2-6) before we
#    touch anything - both the cell values and the hyperlink targets.
# ---------------------------------------------------------------------
$oldRowCount = 5
$oldValues = @()
for ($i = 0; $i -lt $oldRowCount; $i++) {
    $r = 2 + $i
    $oldValues += , @(
        $ws2.Cells.Item($r,1).Value2,
        $ws2.Cells.Item($r,2).Value2,
        $ws2.Cells.Item($r,3).Value2,
        $ws2.Cells.Item($r,4).Value2,
        $ws2.Cells.Item($r,5).Value2,
        $ws2.Cells.Item($r,6).Value2
    )
}

$oldLinks = @()
foreach ($h in $ws2.Hyperlinks) {
    $oldLinks += $h.Address
}

# ---------------------------------------------------------------------
# 2. Append those 5 rows to the bottom of "Previously added" (rows
#    182-186), preserving the row formatting/styles and re-creating the
#    hyperlinks.
# ---------------------------------------------------------------------
$destStart = 182
for ($i = 0; $i -lt $oldRowCount; $i++) {
    $row = $destStart + $i

    $ws1.Range("A181:F181").Copy()
    $ws1.Range("A$row`:F$row").PasteSpecial($xlPasteFormats)

    $vals = $oldValues[$i]
    $ws1.Cells.Item($row,1).Value = $vals[0]
    $ws1.Cells.Item($row,2).Value = $vals[1]
    $ws1.Cells.Item($row,3).Value = $vals[2]
    $ws1.Cells.Item($row,4).Value = $vals[3]
    $ws1.Cells.Item($row,5).Value = $vals[4]
    $ws1.Cells.Item($row,6).Value = $vals[5]

    $ws1.Hyperlinks.Add($ws1.Cells.Item($row,1), $oldLinks[$i])

    # Hyperlinks.Add re-styles column A with the built-in "Hyperlink"
    # style - put the original (non-hyperlink-styled) look back.
    $ws1.Range("A181").Copy()
    $ws1.Cells.Item($row,1).PasteSpecial($xlPasteFormats)
}
$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Replace the content of "New": drop all 5 old hyperlinks/rows and
#    write the 4 newly scraped listings in their place.
# ---------------------------------------------------------------------
$ws2.Range("A2:A6").Hyperlinks.Delete()
$ws2.Rows.Item(6).Delete($xlShiftUp)

$newRows = @(
    @("https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/malinovas-pag/mknjc.html", "18 500 €", "Daugavpils un raj.", "3 ha.",  "44700010164", 45946.60208333333),
    @("https://www.ss.com/msg/lv/real-estate/wood/dobele-and-reg/berzes-pag/kkjpg.html",         "1 400 €",  "Dobele un raj.",      "1 ha.",  "46520030087", 45946.40972222222),
    @("https://www.ss.com/msg/lv/real-estate/wood/jelgava-and-reg/kalnciems/cfghfd.html",        "89 500 €", "Jelgava un raj.",     "18 ha.", "54310030137", 45946.49097222222),
    @("https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/blomes-pag/lkijb.html",          "123 456 €","Valka un raj.",       "6 ha.",  "94460010165", 45945.875)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row  = 2 + $i
    $vals = $newRows[$i]

    $ws2.Cells.Item($row,1).Value = $vals[0]
    $ws2.Cells.Item($row,2).Value = $vals[1]
    $ws2.Cells.Item($row,3).Value = $vals[2]
    $ws2.Cells.Item($row,4).Value = $vals[3]
    $ws2.Cells.Item($row,5).Value = $vals[4]
    $ws2.Cells.Item($row,6).Value = $vals[5]

    $ws2.Hyperlinks.Add($ws2.Cells.Item($row,1), $vals[0])
}
$ws2.Application.CutCopyMode = $false

Write-Output "forests data updated: +5 rows on 'Previously added', 'New' refreshed with 4 rows"
